# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates metrics on row 3 (metrics_sim_with_priors.json) of the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.6216216216216216
$ws.Range("D3").Value = 0.9459459459459459
$ws.Range("H3").Value = 0.6848072562358276
$ws.Range("I3").Value = 0.08997281304973613
$ws.Range("J3").Value = 0.5135135135135135
$ws.Range("K3").Value = 95.02702702702703

$ws.Range("Q3").Value = 9
$ws.Range("R3").Value = 23
$ws.Range("S3").Value = 43
$ws.Range("T3").Value = 102
$ws.Range("U3").Value = 159
$ws.Range("V3").Value = 836
$ws.Range("W3").Value = 822
$ws.Range("X3").Value = 802
$ws.Range("Y3").Value = 743
$ws.Range("Z3").Value = 686

$ws.Range("AF3").Value = 0.989349
$ws.Range("AG3").Value = 0.972781
$ws.Range("AH3").Value = 0.949112
$ws.Range("AI3").Value = 0.87929
$ws.Range("AJ3").Value = 0.8118339999999999
